# edit.ps1 - applies the simmer 4.2.1 blog-post document changes:
#  1. Collapse "But as [Tom Lawton pointed out], until now" -> "But, until now"
#     (removes the hyperlinked run pointing at the Google Groups thread).
#  2. Remove the "Finally, the readership may find interesting ..." paragraph
#     (and its two hyperlinks) entirely, leaving an empty paragraph that takes
#     on the heading ("New features:") paragraph formatting, and add a
#     lastRenderedPageBreak marker to the "New features:" run.
#  3. Remove the trailing "Article originally published in Enchufa2.es: ..."
#     paragraph (and its hyperlink) entirely.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "But as Tom Lawton pointed out, until now" -> "But, until now"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "But as Tom Lawton pointed out, until now",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "But, until now", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: drop the "Finally, ..." paragraph, replace with an empty
# paragraph sharing the "New features:" heading formatting, and stamp the
# page-break marker on the heading run.
# ---------------------------------------------------------------------------
$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$headingPPr = '<w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:outlineLvl w:val="2"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr>'

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Finally, the readership")) {
        $emptyParaXml = '<w:p>' + $headingPPr + '</w:p>'
        $p.Range.InsertXML($xmlHeader + $emptyParaXml + $xmlFooter)
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("New features:")) {
        $headingRunXml = '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:lastRenderedPageBreak/><w:t>New features:</w:t></w:r>'
        $headingParaXml = '<w:p>' + $headingPPr + $headingRunXml + '</w:p>'
        $p.Range.InsertXML($xmlHeader + $headingParaXml + $xmlFooter)
        break
    }
}

# ---------------------------------------------------------------------------
# Change 3: remove the trailing "Article originally published ..." paragraph.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Article originally published")) {
        $p.Range.Delete()
        break
    }
}
